# Apply the edit described by the diff:
# - Add a new raid boss row (row 5) "Frozen Dead King", a duplicate of row 3
#   ("The Ice Queen") except for the name.
# - This introduces a new shared string "Frozen Dead King".
# - Row heights for rows 3 and 4 change from 13.8 to 14.25 (matches sheet's
#   default row height), and the new row 5 also uses 14.25.
# - The worksheet selection moves to B9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 3's ("The Ice Queen") values into the new row 5, changing the name.
$ws.Range("A5").Value = 367
$ws.Range("B5").Value = "Frozen Dead King"
$ws.Range("C5").Value = 16000000000
$ws.Range("D5").Value = 16000000000
$ws.Range("E5").Value = 16000000000
$ws.Range("F5").Value = 16000000000
$ws.Range("G5").Value = 16000000000
$ws.Range("H5").Value = 16000000000
$ws.Range("I5").Value = 16000000000
$ws.Range("J5").Value = 8000000000
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1
$ws.Range("N5").Value = 1
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 1
$ws.Range("S5").Value = 9999
$ws.Range("T5").Value = "int"
$ws.Range("U5").Value = 300000
$ws.Range("V5").Value = 1
$ws.Range("W5").Value = 1000000000000
$ws.Range("X5").Value = 0
$ws.Range("Y5").Value = "5000000000000-10000000000000"
$ws.Range("Z5").Value = "4000000000-8000000000"
$ws.Range("AA5").Value = 500000000
$ws.Range("AB5").Value = 500000000
$ws.Range("AC5").Value = 1
$ws.Range("AD5").Value = 1
$ws.Range("AE5").Value = 1
$ws.Range("AF5").Value = 1
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 1
$ws.Range("AJ5").Value = 1
$ws.Range("AK5").Value = 1
$ws.Range("AL5").Value = 1
$ws.Range("AN5").Value = 0
$ws.Range("AO5").Value = "The Ice Plane"
$ws.Range("AQ5").Value = 0.4
$ws.Range("AR5").Value = 0.75
$ws.Range("AS5").Value = 0.45
$ws.Range("AU5").Value = 1
$ws.Range("AV5").Value = 1
$ws.Range("AW5").Value = 0.35

# Row height adjustments (13.8 -> 14.25 for existing rows 3 & 4; new row 5 is 14.25).
$ws.Rows.Item(3).RowHeight = 14.25
$ws.Rows.Item(4).RowHeight = 14.25
$ws.Rows.Item(5).RowHeight = 14.25

# Update the selected cell to match the post-edit state.
$ws.Range("B9").Select() | Out-Null

Write-Host "Applied raid boss row edit."
